# Adds a new "2023" data column (column S) to the "sastumroebi da
# restornebi" (hotels and restaurants) statistics sheet, carrying the
# Batumi / Ajara figures forward one more year.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Seed column S with the same cell formatting as column R (the previous
# year's column) for every data row, then overwrite with the new values.
$dataRows = 3..14
foreach ($row in $dataRows) {
    $ws.Range("R$row").Copy($ws.Range("S$row"))
}

$ws.Range("S3").Value2  = 2023
$ws.Range("S4").Value2  = 25.6
$ws.Range("S5").Value2  = 25.6
$ws.Range("S6").Value2  = 1006
$ws.Range("S7").Value2  = 971
$ws.Range("S8").Value2  = 916.6
$ws.Range("S9").Value2  = 14.3
$ws.Range("S10").Value2 = 10.7
$ws.Range("S11").Value2 = 11.3
$ws.Range("S12").Value2 = 1.1000000000000001
$ws.Range("S13").Value2 = 13.7
$ws.Range("S14").Value2 = "_"

# Leave the selection where the editor's cursor ended up.
$null = $ws.Range("F18").Select()
